$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3516.6042
$ws.Range("I64").Value = 3410.276
$ws.Range("J64").Value = 3678.8948
$ws.Range("K64").Value = 3410.276
$ws.Range("L64").Value = 3678.8948
$ws.Range("M64").Value = -3162.276
$ws.Range("N64").Value = -4174.8948
# Row 67
$ws.Range("H67").Value = 3516.6042
$ws.Range("I67").Value = 3410.276
$ws.Range("J67").Value = 3678.8948
$ws.Range("K67").Value = 3410.276
$ws.Range("L67").Value = 3678.8948
$ws.Range("M67").Value = -2552.276
$ws.Range("N67").Value = -5394.8948
# Row 116
$ws.Range("H116").Value = 6414916
$ws.Range("I116").Value = 12827280
$ws.Range("J116").Value = 2551
$ws.Range("K116").Value = 12827280
$ws.Range("L116").Value = 2551
$ws.Range("M116").Value = -12823838
$ws.Range("N116").Value = -9435
# Row 132
$ws.Range("H132").Value = 3780.3333
$ws.Range("I132").Value = 3161.923
$ws.Range("J132").Value = 7800
$ws.Range("K132").Value = 9485.769
$ws.Range("L132").Value = 23400
$ws.Range("M132").Value = -6955.769
$ws.Range("N132").Value = -28460
# Row 138
$ws.Range("H138").Value = 1297.88
$ws.Range("I138").Value = 1027.5319
$ws.Range("K138").Value = 3082.5957
$ws.Range("M138").Value = 2057.4043

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2966.7778
$ws.Range("I2").Value = 3183.5
$ws.Range("J2").Value = 2533.3333
$ws.Range("K2").Value = 3183.5
$ws.Range("L2").Value = 2533.3333
$ws.Range("M2").Value = -3070.5
$ws.Range("N2").Value = -2759.3333
# Row 32
$ws.Range("H32").Value = 10993399
$ws.Range("I32").Value = 13516139
$ws.Range("K32").Value = 13516139
$ws.Range("M32").Value = -13515852
# Row 88
$ws.Range("H88").Value = 1184.3334
$ws.Range("I88").Value = 1326.5
$ws.Range("J88").Value = 900
$ws.Range("K88").Value = 1326.5
$ws.Range("L88").Value = 900
$ws.Range("M88").Value = -920.5
$ws.Range("N88").Value = -1712
# Row 91
$ws.Range("H91").Value = 1184.3334
$ws.Range("I91").Value = 1326.5
$ws.Range("J91").Value = 900
$ws.Range("K91").Value = 1326.5
$ws.Range("L91").Value = 900
$ws.Range("M91").Value = 77.5
$ws.Range("N91").Value = -3708
# Row 116
$ws.Range("H116").Value = 2966.7778
$ws.Range("I116").Value = 3183.5
$ws.Range("J116").Value = 2533.3333
$ws.Range("K116").Value = 3183.5
$ws.Range("L116").Value = 2533.3333
$ws.Range("M116").Value = -889.5
$ws.Range("N116").Value = -7121.3333
# Row 122
$ws.Range("H122").Value = 838.80554
$ws.Range("I122").Value = 778.7273
$ws.Range("J122").Value = 1499.6666
$ws.Range("K122").Value = 2336.1819
$ws.Range("L122").Value = 4498.9998
$ws.Range("M122").Value = 113.8181
$ws.Range("N122").Value = -9398.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2966.7778
$ws.Range("I3").Value = 3183.5
$ws.Range("J3").Value = 2533.3333
$ws.Range("K3").Value = 3183.5
$ws.Range("L3").Value = 2533.3333
$ws.Range("M3").Value = -3069.5
$ws.Range("N3").Value = -2761.3333
# Row 20
$ws.Range("H20").Value = 2773.9546
$ws.Range("I20").Value = 1567.7142
$ws.Range("K20").Value = 1567.7142
$ws.Range("M20").Value = -1320.7142
# Row 105
$ws.Range("H105").Value = 2915.9707
$ws.Range("I105").Value = 540
$ws.Range("J105").Value = 3025.6309
$ws.Range("K105").Value = 540
$ws.Range("L105").Value = 3025.6309
$ws.Range("M105").Value = 1207
$ws.Range("N105").Value = -6519.6309

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 450
$ws.Range("I10").Value = 450
$ws.Range("K10").Value = 450
$ws.Range("M10").Value = -311
# Row 94
$ws.Range("H94").Value = 3992
$ws.Range("I94").Value = 1759.4445
$ws.Range("J94").Value = 6001.3
$ws.Range("K94").Value = 1759.4445
$ws.Range("L94").Value = 6001.3
$ws.Range("M94").Value = -1308.4445
$ws.Range("N94").Value = -6903.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 16070519
$ws.Range("I132").Value = 891.8889
$ws.Range("J132").Value = 88383840
$ws.Range("K132").Value = 8027.0001
$ws.Range("L132").Value = 795454560
$ws.Range("M132").Value = -5497.0001
$ws.Range("N132").Value = -795459620

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 728.36365
$ws.Range("I16").Value = 601.6
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 601.6
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -431.6
$ws.Range("N16").Value = -1340
# Row 22
$ws.Range("H22").Value = 677
$ws.Range("I22").Value = 616.25
$ws.Range("J22").Value = 920
$ws.Range("K22").Value = 616.25
$ws.Range("L22").Value = 920
$ws.Range("M22").Value = -321.25
$ws.Range("N22").Value = -1510
# Row 27
$ws.Range("H27").Value = 677
$ws.Range("I27").Value = 616.25
$ws.Range("J27").Value = 920
$ws.Range("K27").Value = 616.25
$ws.Range("L27").Value = 920
$ws.Range("M27").Value = -509.25
$ws.Range("N27").Value = -1134
# Row 64
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20450
# Row 67
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21560
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 15000
$ws.Range("J63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16248
# Row 66
$ws.Range("H66").Value = 15000
$ws.Range("J66").Value = 15000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -51240
# Row 81
$ws.Range("H81").Value = 713.25
$ws.Range("I81").Value = 643.5454999999999
$ws.Range("J81").Value = 1480
$ws.Range("K81").Value = 1287.091
$ws.Range("L81").Value = 2960
$ws.Range("M81").Value = -226.0909999999999
$ws.Range("N81").Value = -5082
# Row 84
$ws.Range("H84").Value = 713.25
$ws.Range("I84").Value = 643.5454999999999
$ws.Range("J84").Value = 1480
$ws.Range("K84").Value = 6435.455
$ws.Range("L84").Value = 14800
$ws.Range("M84").Value = -1131.455
$ws.Range("N84").Value = -25408
